# Update "想去人数" (interested-count) figures on sheet "展览" and sheet "全部类型"
# to reflect the latest scrape output (gh-pages rebuild at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 13385
    $ws.Range("F4").Value = 650
    $ws.Range("F5").Value = 220

    if ($sheetName -eq "展览") {
        $ws.Range("F6").Value = 452
        $ws.Range("F7").Value = 1330
    } elseif ($sheetName -eq "全部类型") {
        $ws.Range("F8").Value = 452
        $ws.Range("F9").Value = 1330
    }
}
